# Documentation dossier de réalisation
# Adds new glossary entries to "Feuil1" and switches the active tab back
# to "Feuil1" (from "Feuil3").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")

# --- New glossary terms appended at the bottom of Feuil1 (column A only) ---
$newTerms = @(
    "Sprint",
    "Merge",
    "Colaborate",
    "Mesh",
    "SVG",
    "Porte logique",
    "ANSI",
    "IEL",
    "DIN",
    "PNG",
    "Input",
    "Output",
    "Line Renderer",
    "Collider"
)

$startRow = 68
for ($i = 0; $i -lt $newTerms.Count; $i++) {
    $row = $startRow + $i
    $ws1.Range("A$row").Value = $newTerms[$i]
}

# --- Switch the active sheet back to Feuil1 (was Feuil3) ---
[void]$ws1.Activate()

# Scroll so row 61 is at the top of the view (best-effort; mirrors topLeftCell="A61")
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1

# Final selection on Feuil1 ends up on B75
[void]$ws1.Range("B75").Select()
